$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.337.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -2.58%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.855.80"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -3.00%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.10%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'326.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.15%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  +0.02%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4554"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -2.70%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3899"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -2.79%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'47.87"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -9.88%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.07928"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -5.06%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'1.012"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -2.76%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'21.44"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -2.96%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.863.64"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.30%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'5.919"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -2.15%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'7.166"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -3.66%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'1.004"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.08%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.06624"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.78%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'86.05"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -3.66%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.00001029"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.98%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'17.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -3.85%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  +0.14%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'5.499"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -3.57%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'27.312.50"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -2.66%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'10.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -3.76%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +0.17%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'2.072.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -1.99%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'153.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.04%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -0.14%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -2.95%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'5.463"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -3.34%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'121.25"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -1.24%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.9483"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -2.27%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.09358"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.98%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.446"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.40%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'3.593"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -1.07%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'5.254"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -4.80%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.06037"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -1.41%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.02227"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -2.94%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'1.210"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.20%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +0.09%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'8.038"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -8.22%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.5920"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -2.96%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.1885"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.59%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'10.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -7.48%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -1.80%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.5617"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -3.55%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'12.05"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -5.62%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -1.93%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.917"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -4.68%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.06739"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.42%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'108.22"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.77%  "
$ws.Range("E51").Style = "Normal"
